$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Unmerge the A2:A3 and B2:B3 merged cells (only A5:C5 stays merged)
$ws.Range("A2:A3").UnMerge()
$ws.Range("B2:B3").UnMerge()

# 2. Update the text content: the "Activity" list now has 3 separate rows
#    instead of one merged "1" row, and the "Result / Actions" column (C)
#    is now empty for each.
$ws.Range("B2").Value = "learned HTTP and web development"
$ws.Range("C2").ClearContents()

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Reviewed the learning notes in the past 2 weeks"
$ws.Range("C3").ClearContents()

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Take the interview by VT department R&D representative"
$ws.Range("C4").ClearContents()

# 3. Row heights: row3 goes back to default, row4 takes the ht=27 that row3 had
$ws.Rows("3").RowHeight = $ws.Rows("1").RowHeight
$ws.Rows("4").RowHeight = 27

# 4. Selection moves to a full-column selection of A:A
$ws.Range("A1:A1048576").Select()

Write-Output "done"
